$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual "B" column (spreadsheet column C) values that toggled
# between blank/imputed across rows 2-23 ---
$ws.Range("C2").Value = 14.9
$ws.Range("C6").ClearContents()
$ws.Range("C12").Value = 12.5
$ws.Range("C14").ClearContents()
$ws.Range("C20").Value = 12.5
$ws.Range("C21").Value = 12.7
$ws.Range("C22").ClearContents()
$ws.Range("C23").ClearContents()

# --- Remove the two rows that were dropped from the data set (RM 232 and
# SC 92), which shifts all subsequent rows up and shrinks the used range
# from A1:F35 to A1:F33 ---
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# --- After the shift, a few more cells changed value (some previously
# blank cells now hold imputed numbers, and vice versa) ---
$ws.Range("B30").Value = -19.7
$ws.Range("C31").Value = 15.3
$ws.Range("B32").ClearContents()
$ws.Range("C33").Value = 10.4
